$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the header row's style (s="4", the "mtitleStyle" cell style) on A10:A12
$ws.Range("A9").Copy()
$ws.Range("A10:A12").PasteSpecial(-4122)

# C11 keeps its current number formatting/style but its text needs to change
# from "-3" to "-1" while staying a text value (not becoming a number). Stash
# the existing format, force Text entry, then restore the original format.
$ws.Range("C11").Copy()
$ws.Range("Z1").PasteSpecial(-4122)
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = "-1"
$ws.Range("Z1").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

# C12 / E12 are plain value updates
$ws.Range("C12").Value = -2
$ws.Range("E12").Value = "88/140"
